$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (copy style from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$timestamps = @(
    "2021-10-05 13:41:06.011840",
    "2021-10-05 13:41:06.011852",
    "2021-10-05 13:41:06.011856",
    "2021-10-05 13:41:06.011859",
    "2021-10-05 13:41:06.011862",
    "2021-10-05 13:41:06.011866",
    "2021-10-05 13:41:06.011869",
    "2021-10-05 13:41:06.011872",
    "2021-10-05 13:41:06.011875",
    "2021-10-05 13:41:06.011878",
    "2021-10-05 13:41:06.011881",
    "2021-10-05 13:41:06.011884",
    "2021-10-05 13:41:06.011887",
    "2021-10-05 13:41:06.011890",
    "2021-10-05 13:41:06.011893",
    "2021-10-05 13:41:06.011896"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
